$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.057860016822815
$ws.Range("B1").Value = 1.348351359367371
$ws.Range("C1").Value = 1.984971880912781
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 2.136033535003662
